{"js": "// Remove the blog-footer boilerplate paragraphs that used to follow the\n// \"LOQ4201: Introdu\u00e7\u00e3o \u00e0 Engenharia de Produ\u00e7\u00e3o (Requisito fraco)\" line:\n//   - one empty paragraph\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//     pages. Original theme under Creative Commons Attribution\"\n// The paragraph that follows them (another empty paragraph, then the\n// page-break paragraph) must remain untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst anchorText = \"LOQ4201: Introdu\u00e7\u00e3o \u00e0 Engenharia de Produ\u00e7\u00e3o (Requisito fraco)\";\nconst footerText1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst footerText2 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // Walk forward from the anchor paragraph, collecting the boilerplate\n  // block to delete: an optional blank paragraph, then the two known\n  // footer-text paragraphs (in order). Stop as soon as the pattern breaks.\n  const candidates = [];\n  let idx = anchorIndex + 1;\n  let blankIdx = -1;\n\n  if (idx < items.length && items[idx].text === \"\") {\n    blankIdx = idx;\n    idx++;\n  }\n  if (idx < items.length && items[idx].text === footerText1) {\n    candidates.push(items[idx]);\n    idx++;\n  }\n  if (idx < items.length && items[idx].text === footerText2) {\n    candidates.push(items[idx]);\n    idx++;\n  }\n\n  // Only remove the leading blank paragraph if at least one footer\n  // paragraph was actually found right after it (otherwise the document\n  // has already been cleaned up and there is nothing to do).\n  const toDelete = [];\n  if (candidates.length > 0 && blankIdx !== -1) {\n    toDelete.push(items[blankIdx]);\n  }\n  toDelete.push(...candidates);\n\n  for (const para of toDelete) {\n    para.delete();\n  }\n\n  if (toDelete.length > 0) {\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the blog-footer boilerplate paragraphs that used to follow the\n# \"LOQ4201: Introdu\u00e7\u00e3o \u00e0 Engenharia de Produ\u00e7\u00e3o (Requisito fraco)\" line:\n#   - one empty paragraph\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#     pages. Original theme under Creative Commons Attribution\"\n# The paragraph that follows them (another empty paragraph, then the\n# page-break paragraph) must remain untouched.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOQ4201: Introdu\u00e7\u00e3o \u00e0 Engenharia de Produ\u00e7\u00e3o (Requisito fraco)\"\n$footerText1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$footerText2 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`n\", \"`v\")\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    # Figure out exactly how many of the following paragraphs are the\n    # boilerplate block to remove: an optional blank paragraph followed by\n    # the two known footer-text paragraphs.\n    $hasBlank = $false\n    $footerCount = 0\n    $idx = $anchorIndex + 1\n\n    if ($idx -le $count) {\n        $t = $d.Paragraphs.Item($idx).Range.Text.TrimEnd(\"`r\", \"`n\", \"`v\")\n        if ($t -eq \"\") {\n            $hasBlank = $true\n            $idx++\n        }\n    }\n\n    if ($idx -le $count) {\n        $t = $d.Paragraphs.Item($idx).Range.Text.TrimEnd(\"`r\", \"`n\", \"`v\")\n        if ($t -eq $footerText1) {\n            $footerCount++\n            $idx++\n        }\n    }\n\n    if ($idx -le $count) {\n        $t = $d.Paragraphs.Item($idx).Range.Text.TrimEnd(\"`r\", \"`n\", \"`v\")\n        if ($t -eq $footerText2) {\n            $footerCount++\n            $idx++\n        }\n    }\n\n    # Only remove the leading blank paragraph if at least one footer\n    # paragraph was actually found right after it (otherwise the document\n    # has already been cleaned up and there is nothing to do).\n    $deleteCount = $footerCount\n    if ($footerCount -gt 0 -and $hasBlank) {\n        $deleteCount++\n    }\n\n    if ($deleteCount -gt 0) {\n        $startPara = $d.Paragraphs.Item($anchorIndex + 1)\n        $endPara = $d.Paragraphs.Item($anchorIndex + $deleteCount)\n        $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n        $deleteRange.Delete()\n    }\n}\n"}
